$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (M2:T2)
$ws.Range("M2").Value = 70.46291600000001
$ws.Range("N2").Value = 211.388748
$ws.Range("O2").Value = 0.5276750397950939
$ws.Range("P2").Value = 0.5276750397950939
$ws.Range("Q2").Value = 17.30158183283334
$ws.Range("R2").Value = 155.7142364955
$ws.Range("S2").Value = 0.5276750397950939
$ws.Range("T2").Value = 0.5276750397950939

# Row 3 (O3, P3, S3, T3 only)
$ws.Range("O3").Value = 0.07361176802536967
$ws.Range("P3").Value = 0.07361176802536967
$ws.Range("S3").Value = 0.07361176802536967
$ws.Range("T3").Value = 0.07361176802536967

# Row 4 (M4:T4)
$ws.Range("M4").Value = 42.505498
$ws.Range("N4").Value = 127.516494
$ws.Range("O4").Value = 0.3183105613832428
$ws.Range("P4").Value = 0.3183105613832428
$ws.Range("Q4").Value = 10.43687082141667
$ws.Range("R4").Value = 93.93183739275001
$ws.Range("S4").Value = 0.3183105613832428
$ws.Range("T4").Value = 0.3183105613832428

# Row 5 (M5:T5)
$ws.Range("M5").Value = 10.73653933333333
$ws.Range("N5").Value = 32.209618
$ws.Range("O5").Value = 0.08040263079629371
$ws.Range("P5").Value = 0.08040263079629371
$ws.Range("Q5").Value = 2.636267762138889
$ws.Range("R5").Value = 23.72640985925
$ws.Range("S5").Value = 0.08040263079629371
$ws.Range("T5").Value = 0.08040263079629371
